$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.828.62"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.00"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7751"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.01"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07350"
$ws.Range("E9").Value = "  +4.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.25"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08140"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7645"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.451"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.896.15"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.91"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.189"
$ws.Range("E16").Value = "  +4.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.795.67"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.90"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.25"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007849"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.138"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.112.81"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.409"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.10"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.039"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.450"
$ws.Range("E30").Value = "  +5.91%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.470"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05566"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.069"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.244"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7533"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9958"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.634"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01926"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.774"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.148.54"
$ws.Range("E41").Value = "  +11.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4440"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.61"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.947"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8516"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.897"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.75"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.105"
$ws.Range("E49").Value = "  +5.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.817"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.487"
$ws.Range("E51").Value = "  +0.52%  "
